$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-21 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-22 Sunday", 2) | Out-Null
$d.Content.Find.Execute("41×40=1640", $true, $false, $false, $false, $false, $true, 1, $false, "41×32=1312", 2) | Out-Null
$d.Content.Find.Execute("13×17=221", $true, $false, $false, $false, $false, $true, 1, $false, "32×81=2592", 2) | Out-Null
$d.Content.Find.Execute("77×72=5544", $true, $false, $false, $false, $false, $true, 1, $false, "36×33=1188", 2) | Out-Null
$d.Content.Find.Execute("19×34=646", $true, $false, $false, $false, $false, $true, 1, $false, "91×71=6461", 2) | Out-Null
$d.Content.Find.Execute("62×86=5332", $true, $false, $false, $false, $false, $true, 1, $false, "44×95=4180", 2) | Out-Null
$d.Content.Find.Execute("77×79=6083", $true, $false, $false, $false, $false, $true, 1, $false, "66×63=4158", 2) | Out-Null
$d.Content.Find.Execute("63×91=5733", $true, $false, $false, $false, $false, $true, 1, $false, "56×26=1456", 2) | Out-Null
$d.Content.Find.Execute("85×62=5270", $true, $false, $false, $false, $false, $true, 1, $false, "91×97=8827", 2) | Out-Null
$d.Content.Find.Execute("83×75=6225", $true, $false, $false, $false, $false, $true, 1, $false, "26×34=884", 2) | Out-Null
$d.Content.Find.Execute("87×59=5133", $true, $false, $false, $false, $false, $true, 1, $false, "53×37=1961", 2) | Out-Null
$d.Content.Find.Execute("81×15=1215", $true, $false, $false, $false, $false, $true, 1, $false, "61×13=793", 2) | Out-Null
$d.Content.Find.Execute("35×72=2520", $true, $false, $false, $false, $false, $true, 1, $false, "79×44=3476", 2) | Out-Null
$d.Content.Find.Execute("59×55=3245", $true, $false, $false, $false, $false, $true, 1, $false, "16×85=1360", 2) | Out-Null
$d.Content.Find.Execute("62×46=2852", $true, $false, $false, $false, $false, $true, 1, $false, "99×95=9405", 2) | Out-Null
$d.Content.Find.Execute("50×34=1700", $true, $false, $false, $false, $false, $true, 1, $false, "27×81=2187", 2) | Out-Null
$d.Content.Find.Execute("31×75=2325", $true, $false, $false, $false, $false, $true, 1, $false, "72×53=3816", 2) | Out-Null
$d.Content.Find.Execute("41×25=1025", $true, $false, $false, $false, $false, $true, 1, $false, "13×89=1157", 2) | Out-Null
$d.Content.Find.Execute("41×71=2911", $true, $false, $false, $false, $false, $true, 1, $false, "83×89=7387", 2) | Out-Null
$d.Content.Find.Execute("41×11=451", $true, $false, $false, $false, $false, $true, 1, $false, "84×91=7644", 2) | Out-Null
$d.Content.Find.Execute("79×52=4108", $true, $false, $false, $false, $false, $true, 1, $false, "78×34=2652", 2) | Out-Null
$d.Content.Find.Execute("90×36=3240", $true, $false, $false, $false, $false, $true, 1, $false, "11×49=539", 2) | Out-Null
$d.Content.Find.Execute("20×33=660", $true, $false, $false, $false, $false, $true, 1, $false, "81×31=2511", 2) | Out-Null
$d.Content.Find.Execute("90×41=3690", $true, $false, $false, $false, $false, $true, 1, $false, "53×60=3180", 2) | Out-Null
$d.Content.Find.Execute("45×53=2385", $true, $false, $false, $false, $false, $true, 1, $false, "17×85=1445", 2) | Out-Null
$d.Content.Find.Execute("88×63=5544", $true, $false, $false, $false, $false, $true, 1, $false, "30×65=1950", 2) | Out-Null
